$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 198.4
$ws.Range("I9").Value = 77.5
$ws.Range("J9").Value = 279
$ws.Range("K9").Value = 77.5
$ws.Range("L9").Value = 279
$ws.Range("M9").Value = 91.5
$ws.Range("N9").Value = -617

$ws.Range("H141").Value = 3732.4
$ws.Range("I141").Value = 3268.3076
$ws.Range("K141").Value = 9804.9228
$ws.Range("M141").Value = -4624.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 5000
$ws.Range("J9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("N9").Value = -5340

$ws.Range("H20").Value = 5000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5540

$ws.Range("H32").Value = 280551.66
$ws.Range("I32").Value = 1065.3667
$ws.Range("J32").Value = 1677983.1
$ws.Range("K32").Value = 1065.3667
$ws.Range("L32").Value = 1677983.1
$ws.Range("M32").Value = -778.3667
$ws.Range("N32").Value = -1678557.1

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("L33").ClearContents()

$ws.Range("H74").Value = 4125
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4126

$ws.Range("H77").Value = 4125
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20632

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -70608

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H99").Value = 1426.1052
$ws.Range("I99").Value = 968.5
$ws.Range("K99").Value = 968.5
$ws.Range("M99").Value = 529.5

$ws.Range("H134").Value = 3002.8
$ws.Range("I134").Value = 3002.8
$ws.Range("K134").Value = 9008.400000000001
$ws.Range("M134").Value = -6473.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 194
$ws.Range("I5").Value = 194
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 582
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = -470
$ws.Range("M5").ClearContents()

$ws.Range("H6").Value = 85.888885
$ws.Range("I6").Value = 67.166664
$ws.Range("K6").Value = 201.499992
$ws.Range("M6").Value = -88.49999199999999

$ws.Range("H11").Value = 1285.1818
$ws.Range("J11").Value = 1499
$ws.Range("L11").Value = 4497
$ws.Range("N11").Value = -4777

$ws.Range("H13").Value = 525
$ws.Range("I13").Value = 525
$ws.Range("K13").Value = 1575
$ws.Range("M13").Value = -1407

$ws.Range("H16").Value = 81.25
$ws.Range("I16").Value = 152.5
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 457.5
$ws.Range("L16").Value = 30
$ws.Range("M16").Value = -284.5
$ws.Range("N16").Value = -376

$ws.Range("H17").Value = 650
$ws.Range("I17").Value = 583.3333
$ws.Range("J17").Value = 750
$ws.Range("K17").Value = 1749.9999
$ws.Range("L17").Value = 2250
$ws.Range("M17").Value = -1580.9999
$ws.Range("N17").Value = -2588

$ws.Range("H22").Value = 3563.5
$ws.Range("I22").Value = 3002
$ws.Range("J22").Value = 3750.6667
$ws.Range("K22").Value = 9006
$ws.Range("L22").Value = 11252.0001
$ws.Range("M22").Value = -8837
$ws.Range("N22").Value = -11590.0001

$ws.Range("H24").Value = 3550
$ws.Range("I24").Value = 100
$ws.Range("J24").Value = 7000
$ws.Range("K24").Value = 300
$ws.Range("L24").Value = 21000
$ws.Range("M24").Value = -70
$ws.Range("N24").Value = -21460

$ws.Range("H27").Value = 3563.5
$ws.Range("I27").Value = 3002
$ws.Range("J27").Value = 3750.6667
$ws.Range("K27").Value = 9006
$ws.Range("L27").Value = 11252.0001
$ws.Range("M27").Value = -8904
$ws.Range("N27").Value = -11456.0001

$ws.Range("H41").Value = 2.5
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = 3
$ws.Range("K41").Value = 6
$ws.Range("L41").Value = 9
$ws.Range("M41").Value = 332
$ws.Range("N41").Value = -685

$ws.Range("H112").Value = 3233.3333
$ws.Range("I112").Value = 2450
$ws.Range("K112").Value = 7350
$ws.Range("M112").Value = -6242

$ws.Range("H116").Value = 1887
$ws.Range("I116").Value = 1814.5
$ws.Range("J116").Value = 2032
$ws.Range("K116").Value = 5443.5
$ws.Range("L116").Value = 6096
$ws.Range("M116").Value = -2001.5
$ws.Range("N116").Value = -12980

$ws.Range("H135").Value = 194
$ws.Range("I135").Value = 194
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 1746
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = 789
$ws.Range("M135").ClearContents()

$ws.Range("H137").Value = 1510
$ws.Range("I137").Value = 1265
$ws.Range("K137").Value = 3795
$ws.Range("M137").Value = 1305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2431.0527
$ws.Range("I102").Value = 912.7273
$ws.Range("J102").Value = 4518.75
$ws.Range("K102").Value = 912.7273
$ws.Range("L102").Value = 4518.75
$ws.Range("M102").Value = 709.2727
$ws.Range("N102").Value = -7762.75

$ws.Range("H122").Value = 41932.617
$ws.Range("I122").Value = 2761.7058
$ws.Range("K122").Value = 8285.117400000001
$ws.Range("M122").Value = -5835.117400000001

$ws.Range("H132").Value = 4002.3333
$ws.Range("I132").Value = 2600
$ws.Range("K132").Value = 7800
$ws.Range("M132").Value = -5270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 1833.3334
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1833.3334
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = -1645.3334
$ws.Range("M46").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H101").Value = 11054.75
$ws.Range("J101").Value = 11054.75
$ws.Range("L101").Value = 11054.75
$ws.Range("N101").Value = -17544.75

$ws.Range("H105").Value = 36410
$ws.Range("J105").Value = 36410
$ws.Range("L105").Value = 36410
$ws.Range("N105").Value = -43398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("L15").ClearContents()

$ws.Range("H54").Value = 24454.545
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21040

$ws.Range("H69").Value = 22999.5
$ws.Range("J69").Value = 22999.5
$ws.Range("L69").Value = 22999.5
$ws.Range("N69").Value = -24497.5

$ws.Range("H72").Value = 22999.5
$ws.Range("J72").Value = 22999.5
$ws.Range("L72").Value = 68998.5
$ws.Range("N72").Value = -76486.5

$ws.Range("H81").Value = 2998.1667
$ws.Range("I81").Value = 2998.1667
$ws.Range("K81").Value = 5996.3334
$ws.Range("M81").Value = -4935.3334

$ws.Range("H84").Value = 2998.1667
$ws.Range("I84").Value = 2998.1667
$ws.Range("K84").Value = 29981.667
$ws.Range("M84").Value = -24677.667

$ws.Range("H103").Value = 13908.223
$ws.Range("J103").Value = 13908.223
$ws.Range("L103").Value = 13908.223
$ws.Range("N103").Value = -16252.223

$ws.Range("H113").Value = 363.57144
$ws.Range("J113").Value = 666.3333
$ws.Range("L113").Value = 1998.9999
$ws.Range("N113").Value = -6338.9999

$ws.Range("H122").Value = 1174.6471
$ws.Range("I122").Value = 1174.6471
$ws.Range("K122").Value = 3523.9413
$ws.Range("M122").Value = -1073.9413

$ws.Range("H132").Value = 2013.0769
$ws.Range("I132").Value = 2013.0769
$ws.Range("K132").Value = 6039.2307
$ws.Range("M132").Value = -3509.2307

$ws.Range("H136").Value = 10499
$ws.Range("J136").Value = 1000
$ws.Range("L136").Value = 3000
$ws.Range("N136").Value = -8100
